$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column B as Text so the numeric-looking membership counts are
# stored as literal text (shared strings) instead of being inferred as Numbers.
$ws.Range("B2:B6").NumberFormat = "@"

# Column B
$ws.Range('B2').Value = '500'
$ws.Range('B3').Value = '600'
$ws.Range('B4').Value = '256'
$ws.Range('B5').Value = '131'
$ws.Range('B6').Value = '400'

# Column C
$ws.Range('C2').Value = 'Yes, FLASCO encompasses community sites. FLASCO''s membership includes both academic and community oncologists, allowing for representation of all oncology practices in Florida.'
$ws.Range('C3').Value = 'Yes, GASCO encompasses community sites. Community sites are an integral part of clinical oncology practice and GASCO aims to represent and support all oncology professionals within the state of Georgia, regardless of practice setting.'
$ws.Range('C4').Value = 'Yes, IOS encompasses community sites. Many community sites are members of the Indiana Oncology Society and participate in their activities and events.'
$ws.Range('C5').Value = 'Yes, because it includes oncology practices in Iowa communities.'
$ws.Range('C6').Value = 'No, MOASC does not encompass community sites. MOASC focuses on medical oncology practices in Southern California.'

# Column D
$ws.Range('D2').Value = 'No, FLASCO is not a policy advocacy organization, They focus more on education, research, and professional development in the field of clinical oncology.'
$ws.Range('D3').Value = 'No, GASCO is not influential on state or local policy. GASCO is a professional organization focused on supporting clinical oncologists and advancing cancer care in Georgia, but it does not have a direct impact on policymaking at the state or local level.'
$ws.Range('D4').Value = 'No, IOS does not have a direct impact on state or local policy as it is a professional organization focused on oncology practices and research.'
$ws.Range('D5').Value = 'No, lack of publicly available information on their influence on policy.'
$ws.Range('D6').Value = 'No, MOASC primarily focuses on providing education and support to medical professionals in Southern California, rather than influencing policy.'

# Column E
$ws.Range('E2').Value = 'Yes, FLASCO provides engagement opportunity with leadership. FLASCO offers networking and leadership development programs for its members to connect with oncology leaders in the state of Florida.'
$ws.Range('E3').Value = 'Yes, GASCO provides engagement opportunities with leadership through various events and initiatives.'
$ws.Range('E4').Value = 'Yes, IOS provides engagement opportunities with leadership. IOS offers networking events, conferences, and collaborations with industry leaders to foster communication and collaboration among members.'
$ws.Range('E5').Value = 'Yes, The IOWA Oncology Society provides engagement opportunities with leadership through networking events and mentorship programs.'
$ws.Range('E6').Value = 'Yes, MOASC provides engagement opportunity with leadership through various networking events, conferences, and educational programs.'

# Column F
$ws.Range('F2').Value = 'No, FLASCO does not provide support for clinical trial recruitment. FLASCO focuses on education, advocacy, and professional development for oncology professionals.'
$ws.Range('F3').Value = 'Yes, GASCO does provide support for clinical trial recruitment. GASCO offers resources and education to assist oncology professionals in recruiting patients for clinical trials.'
$ws.Range('F4').Value = 'No, IOS does not provide support for clinical trial recruitment._IOS focuses on education and advocacy for oncology professionals, not patient recruitment for clinical trials.'
$ws.Range('F5').Value = 'No, IOWA Oncology Society does not provide support for clinical trial recruitment, as their main focus is on education and networking for oncology professionals.'
$ws.Range('F6').Value = 'No, MOASC does not provide support for clinical trial recruitment. They focus on education and advocacy for medical oncology professionals.'

# Column G
$ws.Range('G2').Value = 'No, FLASCO does not have engagement opportunities with payors. FLASCO focuses on education and advocacy for clinical oncologists, not payor relations.'
$ws.Range('G3').Value = 'Yes, GASCO does provide engagement opportunities with payors. GASCO collaborates with payors to ensure the best outcomes for oncology patients and to address issues related to reimbursement and coverage.'
$ws.Range('G4').Value = 'No, IOS does not provide engagement opportunities with payors. IOS focuses primarily on oncology education and advocacy for oncologists, rather than direct engagement with payors.'
$ws.Range('G5').Value = 'No, IOWA Oncology Society does not engage with payors. The focus of the society is on clinical practice and education rather than payer negotiations.'
$ws.Range('G6').Value = 'No, MOASC does not provide engagement opportunity with payors. MOASC represents medical oncologists in Southern California and focuses on education, advocacy, and networking within the specialty.'

# Column H
$ws.Range('H2').Value = 'Yes, FLASCO includes area experts on its board, as it is a society composed of oncology professionals and leaders in Florida.'
$ws.Range('H3').Value = 'Yes, GASCO includes area experts on its board because they are a professional organization representing clinical oncologists in Georgia.'
$ws.Range('H4').Value = 'Yes, the Indiana Oncology Society likely includes area experts on its board in order to provide expertise and guidance in the field of oncology.'
$ws.Range('H5').Value = 'No, the IOWA Oncology Society does not include area experts on its board. , It may be possible that the board members have expertise in other areas related to oncology or leadership.'
$ws.Range('H6').Value = 'Yes, board members are leading oncology experts in Southern California.'

# Column I
$ws.Range('I2').Value = 'Yes, FLASCO is involved in therapeutic research collaborations. FLASCO actively participates in research collaborations with various organizations, institutions, and stakeholders in the field of oncology to advance clinical research and improve patient outcomes.'
$ws.Range('I3').Value = 'Yes, GASCO is involved in therapeutic research collaborations. This can be seen through their active involvement in research studies and clinical trials aimed at advancing cancer care.'
$ws.Range('I4').Value = 'Yes, collaboration with other medical organizations is essential for advancing therapeutic research.'
$ws.Range('I5').Value = 'No, IOWA Oncology Society primarily focuses on education and advocacy for oncology professionals in Iowa.'
$ws.Range('I6').Value = 'Yes, MOASC is involved in therapeutic research collaborations. The association actively collaborates with pharmaceutical companies, academic institutions, and other organizations to advance research in medical oncology.'

# Column J
$ws.Range('J2').Value = 'No, FLASCO does not include top therapeutic area experts on its board., The organization''s board primarily consists of oncologists and healthcare administrators.'
$ws.Range('J3').Value = 'No,  GASCO does not include top therapeutic area experts on its board. The board may consist of other professionals or stakeholders in the field of oncology, rather than experts in specific therapeutic areas.'
$ws.Range('J4').Value = 'No, IOS does not include top therapeutic area experts on its board. The organization focuses on providing education and support to healthcare professionals in the field of oncology.'
$ws.Range('J5').Value = 'Yes, the IOWA Oncology Society includes top therapeutic area experts on its board because they are oncology specialists.'
$ws.Range('J6').Value = 'No, MOASC does not include top therapeutic area experts on its board. The focus is more on community oncologists rather than specific therapeutic area experts.'

# Column K
$ws.Range('K2').Value = 'Florida'
$ws.Range('K3').Value = 'Georgia'
$ws.Range('K4').Value = 'Midwest.'
$ws.Range('K5').Value = 'Midwest'
$ws.Range('K6').Value = 'Southern California'

# Column B only needed the Text number-format to force string storage; restore
# the default (unstyled) appearance now that the values are text.
$ws.Range("B2:B6").Style = "Normal"
